# Daily attendance processing - 2025-12-31 20:34:28
#
# The "Recorded By" column (G) lists the users who recorded/updated a
# session's attendance. Swap the order of the two names so the human
# editor ("dnasr281@gmail.com") is listed before "System" wherever both
# appear together (i.e. turn "System, dnasr281@gmail.com" into
# "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
